$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H8").Value = 94.8
$ws_ALC.Range("I8").Value = 94.8
$ws_ALC.Range("J8").Value = 0
$ws_ALC.Range("K8").Value = 284.4
$ws_ALC.Range("L8").Value = 0
$ws_ALC.Range("M8").ClearContents()
$ws_ALC.Range("N8").ClearContents()

$ws_ALC.Range("H32").Value = 5000
$ws_ALC.Range("I32").Value = 5000
$ws_ALC.Range("J32").Value = 0
$ws_ALC.Range("K32").Value = 5000
$ws_ALC.Range("L32").Value = 0
$ws_ALC.Range("M32").Value = -4674

$ws_ALC.Range("H98").Value = 5356.4375
$ws_ALC.Range("I98").Value = 3046.8667
$ws_ALC.Range("J98").Value = 40000
$ws_ALC.Range("K98").Value = 3046.8667
$ws_ALC.Range("L98").Value = 40000
$ws_ALC.Range("M98").Value = -1548.8667
$ws_ALC.Range("N98").Value = -42996

$ws_ALC.Range("H122").Value = 5356.4375
$ws_ALC.Range("I122").Value = 3046.8667
$ws_ALC.Range("J122").Value = 40000
$ws_ALC.Range("K122").Value = 9140.6001
$ws_ALC.Range("L122").Value = 120000
$ws_ALC.Range("M122").Value = -6690.6001
$ws_ALC.Range("N122").Value = -124900

$ws_ALC.Range("H137").Value = 1230.8
$ws_ALC.Range("I137").Value = 1230.8
$ws_ALC.Range("J137").Value = 0
$ws_ALC.Range("K137").Value = 3692.4
$ws_ALC.Range("L137").Value = 0
$ws_ALC.Range("M137").Value = -1142.4

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 5359.1113
$ws_ARM.Range("I45").Value = 1748.8572
$ws_ARM.Range("J45").Value = 17995
$ws_ARM.Range("K45").Value = 1748.8572
$ws_ARM.Range("L45").Value = 17995
$ws_ARM.Range("M45").Value = -1371.8572
$ws_ARM.Range("N45").Value = -18749

$ws_ARM.Range("H102").Value = 2499.5
$ws_ARM.Range("I102").Value = 1999
$ws_ARM.Range("J102").Value = 3000
$ws_ARM.Range("K102").Value = 1999
$ws_ARM.Range("L102").Value = 3000
$ws_ARM.Range("M102").Value = -377
$ws_ARM.Range("N102").Value = -6244

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 642.5
$ws_BSM.Range("I20").Value = 642.5
$ws_BSM.Range("J20").Value = 0
$ws_BSM.Range("K20").Value = 642.5
$ws_BSM.Range("L20").Value = 0
$ws_BSM.Range("M20").Value = -395.5

$ws_BSM.Range("H107").Value = 1700
$ws_BSM.Range("I107").Value = 1700
$ws_BSM.Range("J107").Value = 0
$ws_BSM.Range("K107").Value = 1700
$ws_BSM.Range("L107").Value = 0
$ws_BSM.Range("M107").Value = 220

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 499.5
$ws_CRP.Range("I22").Value = 0
$ws_CRP.Range("J22").Value = 499.5
$ws_CRP.Range("K22").Value = 0
$ws_CRP.Range("L22").Value = 499.5
$ws_CRP.Range("N22").Value = -1199.5

$ws_CRP.Range("H33").Value = 42366
$ws_CRP.Range("I33").Value = 1281.25
$ws_CRP.Range("J33").Value = 58799.9
$ws_CRP.Range("K33").Value = 1281.25
$ws_CRP.Range("L33").Value = 58799.9
$ws_CRP.Range("M33").Value = -902.25
$ws_CRP.Range("N33").Value = -59557.9

$ws_CRP.Range("H38").Value = 17602.8
$ws_CRP.Range("I38").Value = 4503.5
$ws_CRP.Range("J38").Value = 70000
$ws_CRP.Range("K38").Value = 4503.5
$ws_CRP.Range("L38").Value = 70000
$ws_CRP.Range("M38").Value = -4126.5
$ws_CRP.Range("N38").Value = -70754

$ws_CRP.Range("H46").Value = 17602.8
$ws_CRP.Range("I46").Value = 4503.5
$ws_CRP.Range("J46").Value = 70000
$ws_CRP.Range("K46").Value = 4503.5
$ws_CRP.Range("L46").Value = 70000
$ws_CRP.Range("M46").Value = -4292.5
$ws_CRP.Range("N46").Value = -70422

$ws_CRP.Range("H88").Value = 32792.625
$ws_CRP.Range("I88").Value = 0
$ws_CRP.Range("J88").Value = 32792.625
$ws_CRP.Range("K88").Value = 0
$ws_CRP.Range("L88").Value = 32792.625
$ws_CRP.Range("N88").Value = -33604.625

$ws_CRP.Range("H91").Value = 32792.625
$ws_CRP.Range("I91").Value = 0
$ws_CRP.Range("J91").Value = 32792.625
$ws_CRP.Range("K91").Value = 0
$ws_CRP.Range("L91").Value = 32792.625
$ws_CRP.Range("N91").Value = -35600.625

$ws_CRP.Range("H106").Value = 77660.11
$ws_CRP.Range("I106").Value = 0
$ws_CRP.Range("J106").Value = 77660.11
$ws_CRP.Range("K106").Value = 0
$ws_CRP.Range("L106").Value = 77660.11
$ws_CRP.Range("N106").Value = -80184.11

$ws_CRP.Range("H107").Value = 609.9
$ws_CRP.Range("I107").Value = 450
$ws_CRP.Range("J107").Value = 1249.5
$ws_CRP.Range("K107").Value = 450
$ws_CRP.Range("L107").Value = 1249.5
$ws_CRP.Range("M107").Value = 1470
$ws_CRP.Range("N107").Value = -5089.5

$ws_CRP.Range("H121").Value = 0
$ws_CRP.Range("I121").Value = 0
$ws_CRP.Range("J121").Value = 0
$ws_CRP.Range("K121").Value = 0
$ws_CRP.Range("L121").ClearContents()
$ws_CRP.Range("N121").ClearContents()

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 1750.2273
$ws_CUL.Range("I4").Value = 1395
$ws_CUL.Range("J4").Value = 4000
$ws_CUL.Range("K4").Value = 4185
$ws_CUL.Range("L4").Value = 12000
$ws_CUL.Range("M4").Value = -4073
$ws_CUL.Range("N4").Value = -12224

$ws_CUL.Range("H20").Value = 0
$ws_CUL.Range("I20").Value = 0
$ws_CUL.Range("J20").Value = 0
$ws_CUL.Range("K20").Value = 0
$ws_CUL.Range("L20").Value = 0
$ws_CUL.Range("M20").ClearContents()

$ws_CUL.Range("H22").Value = 3126.75
$ws_CUL.Range("I22").Value = 0
$ws_CUL.Range("J22").Value = 3126.75
$ws_CUL.Range("K22").Value = 0
$ws_CUL.Range("L22").Value = 9380.25
$ws_CUL.Range("N22").Value = -9718.25

$ws_CUL.Range("H27").Value = 3126.75
$ws_CUL.Range("I27").Value = 0
$ws_CUL.Range("J27").Value = 3126.75
$ws_CUL.Range("K27").Value = 0
$ws_CUL.Range("L27").Value = 9380.25
$ws_CUL.Range("N27").Value = -9584.25

$ws_CUL.Range("H122").Value = 0
$ws_CUL.Range("I122").Value = 0
$ws_CUL.Range("J122").Value = 0
$ws_CUL.Range("K122").Value = 0
$ws_CUL.Range("L122").ClearContents()
$ws_CUL.Range("M122").ClearContents()
$ws_CUL.Range("N122").ClearContents()

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H7").Value = 5000000
$ws_GSM.Range("I7").Value = 5000000
$ws_GSM.Range("J7").Value = 0
$ws_GSM.Range("K7").Value = 5000000
$ws_GSM.Range("L7").Value = 0
$ws_GSM.Range("M7").Value = -4999888

$ws_GSM.Range("H8").Value = 5000000
$ws_GSM.Range("I8").Value = 5000000
$ws_GSM.Range("J8").Value = 0
$ws_GSM.Range("K8").Value = 5000000
$ws_GSM.Range("L8").Value = 0
$ws_GSM.Range("M8").Value = -4999861

$ws_GSM.Range("H70").Value = 0
$ws_GSM.Range("I70").Value = 0
$ws_GSM.Range("J70").Value = 0
$ws_GSM.Range("K70").Value = 0
$ws_GSM.Range("L70").Value = 0
$ws_GSM.Range("M70").ClearContents()

$ws_GSM.Range("H73").Value = 0
$ws_GSM.Range("I73").Value = 0
$ws_GSM.Range("J73").Value = 0
$ws_GSM.Range("K73").Value = 0
$ws_GSM.Range("L73").Value = 0
$ws_GSM.Range("M73").ClearContents()

$ws_GSM.Range("H102").Value = 3908.1765
$ws_GSM.Range("I102").Value = 3840.375
$ws_GSM.Range("J102").Value = 4993
$ws_GSM.Range("K102").Value = 3840.375
$ws_GSM.Range("L102").Value = 4993
$ws_GSM.Range("M102").Value = -2218.375
$ws_GSM.Range("N102").Value = -8237

$ws_GSM.Range("H103").Value = 30000
$ws_GSM.Range("I103").Value = 0
$ws_GSM.Range("J103").Value = 30000
$ws_GSM.Range("K103").Value = 0
$ws_GSM.Range("L103").Value = 30000
$ws_GSM.Range("N103").Value = -32344

$ws_GSM.Range("H123").Value = 79999
$ws_GSM.Range("I123").Value = 0
$ws_GSM.Range("J123").Value = 79999
$ws_GSM.Range("K123").Value = 0
$ws_GSM.Range("L123").Value = 79999
$ws_GSM.Range("N123").Value = -84899

$ws_GSM.Range("H126").Value = 4837.25
$ws_GSM.Range("I126").Value = 4837.25
$ws_GSM.Range("J126").Value = 0
$ws_GSM.Range("K126").Value = 14511.75
$ws_GSM.Range("L126").Value = 0
$ws_GSM.Range("M126").ClearContents()
$ws_GSM.Range("N126").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4832.8
$ws_LTW.Range("I7").Value = 4116
$ws_LTW.Range("J7").Value = 7700
$ws_LTW.Range("K7").Value = 4116
$ws_LTW.Range("L7").Value = 7700
$ws_LTW.Range("M7").Value = -4004
$ws_LTW.Range("N7").Value = -7924

$ws_LTW.Range("H40").Value = 8179.4287
$ws_LTW.Range("I40").Value = 5668
$ws_LTW.Range("J40").Value = 10063
$ws_LTW.Range("K40").Value = 5668
$ws_LTW.Range("L40").Value = 10063
$ws_LTW.Range("M40").Value = -5532
$ws_LTW.Range("N40").Value = -10335

$ws_LTW.Range("H93").Value = 1474.875
$ws_LTW.Range("I93").Value = 1474.875
$ws_LTW.Range("J93").Value = 0
$ws_LTW.Range("K93").Value = 1474.875
$ws_LTW.Range("L93").Value = 0
$ws_LTW.Range("M93").Value = -226.875

$ws_LTW.Range("H126").Value = 4832.8
$ws_LTW.Range("I126").Value = 4116
$ws_LTW.Range("J126").Value = 7700
$ws_LTW.Range("K126").Value = 12348
$ws_LTW.Range("L126").Value = 23100
$ws_LTW.Range("M126").Value = -9878
$ws_LTW.Range("N126").Value = -28040

$ws_LTW.Range("H132").Value = 1700.4
$ws_LTW.Range("I132").Value = 1700.4
$ws_LTW.Range("J132").Value = 0
$ws_LTW.Range("K132").Value = 5101.200000000001
$ws_LTW.Range("L132").Value = 0
$ws_LTW.Range("M132").Value = -2571.200000000001

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H96").Value = 1943.75
$ws_WVR.Range("I96").Value = 1925
$ws_WVR.Range("J96").Value = 2000
$ws_WVR.Range("K96").Value = 1925
$ws_WVR.Range("L96").Value = 2000
$ws_WVR.Range("M96").Value = -552
$ws_WVR.Range("N96").Value = -4746

$ws_WVR.Range("H107").Value = 1366
$ws_WVR.Range("I107").Value = 750
$ws_WVR.Range("J107").Value = 1982
$ws_WVR.Range("K107").Value = 2250
$ws_WVR.Range("L107").Value = 5946
$ws_WVR.Range("M107").Value = -330
$ws_WVR.Range("N107").Value = -9786

$ws_WVR.Range("H126").Value = 1612
$ws_WVR.Range("I126").Value = 1612
$ws_WVR.Range("J126").Value = 0
$ws_WVR.Range("K126").Value = 4836
$ws_WVR.Range("L126").Value = 0
$ws_WVR.Range("M126").Value = -2366

$ws_WVR.Range("H132").Value = 1000.8571
$ws_WVR.Range("I132").Value = 1000.8571
$ws_WVR.Range("J132").Value = 0
$ws_WVR.Range("K132").Value = 3002.5713
$ws_WVR.Range("L132").Value = 0
$ws_WVR.Range("M132").Value = -472.5712999999996

$ws_WVR.Range("H133").Value = 10714.5
$ws_WVR.Range("I133").Value = 0
$ws_WVR.Range("J133").Value = 10714.5
$ws_WVR.Range("K133").Value = 0
$ws_WVR.Range("L133").Value = 10714.5
$ws_WVR.Range("N133").Value = -20834.5

$ws_WVR.Range("H136").Value = 2874.6667
$ws_WVR.Range("I136").Value = 2812
$ws_WVR.Range("J136").Value = 3000
$ws_WVR.Range("K136").Value = 8436
$ws_WVR.Range("L136").Value = 9000
$ws_WVR.Range("M136").Value = -5886
$ws_WVR.Range("N136").Value = -14100
